$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'62.195.75"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.28%  "
$cell = $ws.Range("D3")
$cell.Value = "'2.415.54"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.74%  "
$cell = $ws.Range("D4")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'561.14"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "
$cell = $ws.Range("D6")
$cell.Value = "'143.12"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.80%  "
$cell = $ws.Range("D9")
$cell.Value = "'2.409.38"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("E13").Value = "  +0.10%  "
$cell = $ws.Range("D14")
$cell.Value = "'25.69"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("E15").Value = "  +0.92%  "
$cell = $ws.Range("D16")
$cell.Value = "'2.853.24"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "
$cell = $ws.Range("D17")
$cell.Value = "'62.106.12"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "
$cell = $ws.Range("D18")
$cell.Value = "'2.414.36"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.62%  "
$cell = $ws.Range("D19")
$cell.Value = "'11.27"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D20")
$cell.Value = "'323.33"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D21")
$cell.Value = "'4.17"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  +2.47%  "
$cell = $ws.Range("D24")
$cell.Value = "'65.77"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  -3.00%  "
$cell = $ws.Range("D26")
$cell.Value = "'9.03"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.48%  "
$cell = $ws.Range("D27")
$cell.Value = "'577.43"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +7.87%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$cell = $ws.Range("D28")
$cell.Value = "'2.535.31"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$cell = $ws.Range("D29")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "
$cell = $ws.Range("D30")
$cell.Value = "'0.0₃0945"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.00%  "
$cell = $ws.Range("D31")
$cell.Value = "'8.21"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  +1.36%  "
$cell = $ws.Range("D35")
$cell.Value = "'1.53"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$cell = $ws.Range("D37")
$cell.Value = "'5.61"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$cell = $ws.Range("D38")
$cell.Value = "'4.73"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D39")
$cell.Value = "'153.05"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell = $ws.Range("D40")
$cell.Value = "'0.382"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.61%  "
$cell = $ws.Range("D41")
$cell.Value = "'18.68"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +2.34%  "
$cell = $ws.Range("D45")
$cell.Value = "'148.83"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("E47").Value = "  +1.21%  "
$cell = $ws.Range("D48")
$cell.Value = "'20.10"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$cell = $ws.Range("D49")
$cell.Value = "'0.593"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("E51").Value = "  +1.64%  "
